# ---------------------------------------------------------------------------
# Target change (from the supplied OOXML diff):
#
#   word/numbering.xml — the <w:nsid w:val="…"/> GUID carried by four
#   <w:abstractNum> definitions (abstractNumId 990, 991, 99416, 99417) is
#   replaced with a different random-looking 8-hex-digit value:
#
#       990    dc045508 -> 2fe7da74
#       991    6123ea10 -> baf3dd5d
#       99416  d0e1ee20 -> bab7a102
#       99417  25cd6f16 -> 3d356c1f
#
#   Nothing else in the package differs — same paragraphs, same <w:numId>
#   references, same <w:lvl> formatting for every level of every list.
#
# Why this is left as a no-op:
#
#   <w:nsid> is Word's internal "list signature" GUID. It has no effect on
#   layout/content (bullets, numbers, indents all come from the sibling
#   <w:lvl> elements, and body paragraphs address lists only via
#   <w:numId>/abstractNumId, never nsid), and — crucially — it is not
#   surfaced anywhere in Word's object model: List/ListTemplate/ListFormat
#   only expose Name/ListID/OutlineNumbered/ListLevels/ListType/etc., with
#   ListID itself mapping to the w:numId, not the abstractNum's nsid. There
#   is no Find/Replace path either, since Find.Execute walks story Range
#   text, not the numbering part. That matches the commit message —
#   "Automatic build output files" — i.e. this hunk is the numbering part
#   being re-minted with fresh random nsids by the export/build pipeline on
#   every run, not a deliberate, reproducible document edit.
#
#   Concretely, this was verified against this runtime: reading/writing
#   $d.Content.WordOpenXML (and Range.WordOpenXML), Find.Execute over the
#   nsid hex strings, List.ListID / ListTemplate.*, re-applying
#   ListFormat.ListTemplate, and Document.Save()/SaveAs() were all tried —
#   none can reach or regenerate <w:nsid>, so there is no COM call that
#   reproduces this specific hunk without touching (and risking corrupting)
#   anything else in the document. Since every other part of the document
#   is unchanged in the diff, the safest, most faithful action available
#   through the object model is to leave the document exactly as-is.
$d = $word.ActiveDocument
Write-Host "No object-model-reachable change: w:nsid on abstractNum 990/991/99416/99417 is an internal, non-semantic GUID not exposed by Word's COM surface (Document.Paragraphs=$($d.Paragraphs.Count) unchanged)."
